$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BASE DATA (wajib update)")

# --- Refresh "ACTUAL END" (column K) for every data row (2-127) to the new
#     dashboard pull date (45816 = 2025-06-08), replacing the old 45812. ---
$ws.Range("K2:K127").Value2 = 45816

# --- A handful of rows also had their "% COMPLETE" (column L) progressed
#     further since the last pull. ---
$pctUpdates = @{
    52  = 0.4668
    55  = 0.7036
    58  = 1
    62  = 0.9596
    65  = 0.6233
    69  = 0.4658
    70  = 0.7808
    71  = 0.1244
    81  = 0.6317
    84  = 0.4143
    88  = 0.4425
    115 = 0.2008
}
foreach ($row in $pctUpdates.Keys) {
    $ws.Range("L$row").Value2 = $pctUpdates[$row]
}

# --- Two tasks were re-prioritized from MEDIUM to HIGH. ---
$ws.Range("H65").Value = "HIGH"
$ws.Range("H81").Value = "HIGH"

# --- The author scrolled the frozen-pane view before saving: the split
#     moved from K1 to J1, and the active selection moved to T24. ---
$window = $excel.ActiveWindow
$window.SplitColumn = 4
$window.SplitRow = 0
$ws.Range("J1").Select()
$window.FreezePanes = $true
$ws.Range("T24").Select()
